$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "client code" header to "clientNo"
$ws.Range("D1").Value = "clientNo"

# Add new "groupNo" column header
$ws.Range("E1").Value = "groupNo"

# Move selection to E2, matching the saved selection state
$ws.Range("E2").Select()
